$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new column before column A: existing Categorie/Opdracht/... shift
# one column to the right (B:F); their styles/row-heights move with them.
$ws.Range("A1").EntireColumn.Insert()

# Give the new header cell (A1) the same formatting as the other header
# cells, then set its text.
$ws.Range("B1").Copy()
$ws.Range("A1").PasteSpecial(-4122)
$excel.CutCopyMode = $false
$ws.Range("A1").Value = "Hoofdcategorie"

# Rows 2-4 (the "Anatomie" questions) keep a Hoofdcategorie of "Anatomie" in
# the new column A; match the formatting of their row first.
$ws.Range("B2").Copy()
$ws.Range("A2:A4").PasteSpecial(-4122)
$excel.CutCopyMode = $false
$ws.Range("A2:A4").Value = "Anatomie"

# The old Categorie column (now B) becomes the subcategory for those rows.
$ws.Range("B2").Value = "enkel"
$ws.Range("B3").Value = "voet"
$ws.Range("B4").Value = "knie"

# Rows 5-11 are left with no Hoofdcategorie value at all (column A stays
# blank/untouched for them); the old Categorie text ("Casus Enkelletsel",
# "Gedrag & Communicatie") simply slid into column B by the insert above.

$ws.Range("B5").Select()
